# OO_Bond_Scans.xlsx update
#
# Change the reported non-Hartree energy column from kcal/mol to kJ/mol
# (multiply all values by 4.184), re-label the header, turn the
# "0 K energy barrier" cell (F2) into a live formula of the original
# kcal/mol number times 4.184, and strip the inherited numeric style
# from the recomputed Energy column (it now uses the default style).
# Also refreshes the selected cell / active sheet bookmarks the workbook
# was left with.

$wb = $excel.ActiveWorkbook

$lastRows = @{ 1 = 7; 2 = 7; 3 = 9; 4 = 10; 5 = 12; 6 = 10 }
$selections = @{
    1 = "C1"
    2 = "F8"
    3 = "C2:C9"
    4 = "C2:C10"
    5 = "C2:C12"
    6 = "F2"
}

for ($i = 1; $i -le 6; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $lastRow = $lastRows[$i]

    # Header: "Energy (kcal/mol)" -> "Energy (kJ/mol)"
    $ws.Range("C1").Value = "Energy (kJ/mol)"

    # Convert the barrier cell (if present on this sheet) into a formula
    # expressing the kJ/mol value in terms of the original kcal/mol number.
    $barrierCell = $ws.Cells.Item(2, 6)
    if ($barrierCell.Value2 -ne $null) {
        $oldKcal = $barrierCell.Value2
        $barrierCell.Formula = "=" + $oldKcal + "*4.184"
    }

    # Recompute the Energy column values (kcal/mol -> kJ/mol) and drop the
    # inherited numeric style from those cells.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        $oldVal = $cell.Value2
        if ($oldVal -ne $null) {
            $cell.Value2 = $oldVal * 4.184
            $cell.Style = "Normal"
        }
    }
}

# Restore per-sheet selections.
for ($i = 1; $i -le 6; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Activate()
    $ws.Range($selections[$i]).Select()
}

# Leave "GAlt-Mono-T-0" (sheet 6) as the active tab, matching the saved file.
$wb.Worksheets.Item(6).Activate()
$wb.Worksheets.Item(6).Range("F2").Select()
